# "add dars from after add 1p5-6"
# Adds new "dars" (lesson) 5 and 6 dialogue rows (168-186, 0-based ids)
# to the "kalimat durusul lughoh 2" sheet, un-hides the previously
# hidden rows, updates row 168 (juz/tamrin columns) and repositions the
# frozen pane / selection to the newly added block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill in the two missing values on the existing last row (168) ---
# Before: only C168(dars)=2, E168(hal)=11, F168(kalimat) were populated.
# After : B168(juz)=1 and D168(tamrin)=0 are also populated.
$ws.Range("B168").Value = 1
$ws.Range("D168").Value = 0

# --- 2. Append the new rows (169-187) describing the new dialogue ---
$newRows = @(
    @(169, 168, 2, 1, 0, 5, "bagaimana kabarmu wahai Ustadz ? Semoga kamu dalam keadaan baik"),
    @(170, 169, 2, 1, 0, 5, "alhamdulillah, dan bagaimana kabarmu wahai hasyim, ? Saya mencintaimu dengan banyak wahai hasyim, sesungguhnya engkau adalah murid yang pandai dan rajin dan memiliki akhlaq yang baik. Apakah dari pakistan kamu atau dari india wahai hasyim ?"),
    @(171, 170, 2, 1, 0, 5, "sesungguhnya aku dari india"),
    @(172, 171, 2, 1, 0, 5, "dan temanku yang keluar bersamamu sekarang dari kelas, apakah dia juga dari india ?"),
    @(173, 172, 2, 1, 0, 5, "tidak, sesungguhnya dia dari pakistan"),
    @(174, 173, 2, 1, 0, 5, "sesungguhnya jam tangan mu indah wahai hasyim, apakah dari jepang dia (jam) ?"),
    @(175, 174, 2, 1, 0, 5, "tidak, sesungguhnya dia (jam) dari india"),
    @(176, 175, 2, 1, 0, 5, "apakah dia (jam) mahal atau murah ?"),
    @(177, 176, 2, 1, 0, 5, "sesungguhnya dia (jam) murah sekali. Sesungguhnya dengan 100 rupiah saja"),
    @(178, 177, 2, 1, 0, 5, "berapa saudara laki laki (milik) mu wahai hasyim ?"),
    @(179, 178, 2, 1, 0, 6, "aku memiliki 3 saudara laki-laki"),
    @(180, 179, 2, 1, 0, 6, "apakah murid-murid (mereka) ?"),
    @(181, 180, 2, 1, 0, 6, "tida sesungguhnya mereka (lk) (adalah) para pedagang"),
    @(182, 181, 2, 1, 0, 6, "dan berapa saudara (pr) milikmu ?"),
    @(183, 182, 2, 1, 0, 6, "aku memiliki 4 saudara (pr)"),
    @(184, 183, 2, 1, 0, 6, "apakah di india mereka (pr) sekarang ?"),
    @(185, 184, 2, 1, 0, 6, "tidak, sesungguhnya mereka (pr) dengan (di) madinah munawaroh bersama bapakku dan ibuku"),
    @(186, 185, 2, 1, 0, 6, "apakah murid-murid (mereka pr) ?"),
    @(187, 186, 2, 1, 0, 6, "tidak, sesungguhnya mereka (pr) adalah guru-guru (pr) madrasah tsanawiyah")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
}

# --- 3. Un-hide the rows that were previously collapsed (2-124) ---
$ws.Rows("2:124").Hidden = $false

# --- 4. Re-position the frozen header pane and the active selection ---
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$win.ScrollRow = 177
$ws.Range("A167:A187").Select()
